# Updates the Price (column D) and Volume(1h) (column E) figures for the
# cryptocurrency table to reflect the latest scrape performed by the
# "Updated symbol list" GitHub Actions workflow.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "305.35"
    "E2" = "1.57%"
    "D3" = "35.87"
    "E3" = "1.04%"
    "D4" = "5.052"
    "E4" = "-0.33%"
    "D5" = "0.08098"
    "E5" = "1.25%"
    "D6" = "1.921"
    "E6" = "1.38%"
    "D7" = "4.150"
    "E7" = "2.45%"
    "D8" = "7.843"
    "E8" = "1.13%"
    "D9" = "0.9309"
    "E9" = "0.22%"
    "D10" = "0.1279"
    "E10" = "-11.75%"
    "D11" = "0.1920"
    "E11" = "1.11%"
    "D12" = "0.09238"
    "E12" = "2.93%"
    "D13" = "0.03492"
    "E13" = "0.95%"
    "D14" = "0.09935"
    "E14" = "0.91%"
    "D15" = "0.001418"
    "E15" = "1.75%"
    "E16" = "16.13%"
    "D17" = "3.616"
    "E17" = "2.42%"
    "D19" = "0.3445"
    "E19" = "0.04%"
    "D20" = "5.191"
    "E20" = "2.62%"
    "E21" = "0.08%"
    "E22" = "5.44%"
    "D23" = "0.04407"
    "E23" = "-1.97%"
    "E24" = "1.70%"
    "D25" = "0.004733"
    "E25" = "-0.53%"
    "D27" = "0.0003128"
    "E27" = "3.46%"
    "E39" = "9.16%"
    "D40" = "0.05163"
    "E40" = "9.03%"
    "D41" = "0.007536"
    "E41" = "2.82%"
    "D42" = "0.01008"
    "E42" = "-5.00%"
    "D43" = "0.1373"
    "E43" = "3.51%"
    "E44" = "-0.48%"
    "E45" = "-1.30%"
    "D46" = "0.00006338"
    "E46" = "1.49%"
    "E47" = "0.02%"
    "D48" = "63.57"
    "E48" = "-1.70%"
    "E49" = "-3.58%"
    "D50" = "0.00002100"
    "E50" = "0.02%"
    "E51" = "0.02%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force a text number format so Excel keeps the value as a literal
    # string (matching the inline-string cells already in the sheet)
    # instead of silently coercing it into a Double and losing the
    # original formatting (trailing zeros, percent sign, etc).
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
